$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# Row 6 (resource id 3, "道心"): Effects formula changed from TATTRI,0,0 to USEHEART,0,0
$ws.Range("G6").Value = "USEHEART,0,0"

# Rows 46-69 ("地图..." resources): UsedEffect was a shared "OPENEVENT,0,0" placeholder
# for every map; it is now a unique OPENEVENT,<id>,0 per row, and the old highlight
# formatting (red/yellow fill used while the value was still a TODO placeholder) is
# cleared now that every row has its own real event id.
$openEventStart = 119
for ($row = 46; $row -le 69; $row++) {
    $eventId = $openEventStart + ($row - 46)
    $cell = $ws.Range("I$row")
    $cell.Value = "OPENEVENT,$eventId,0"
    $cell.ClearFormats()
}

# --- View / layout cosmetics ------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 18.714285714285715
$ws.Columns.Item(7).ColumnWidth = 13.428571428571429
$ws.Columns.Item(9).ColumnWidth = 13.714285714285714

$ws.Range("A22").Select()
$ws.Range("I29").Select()

$excel.ActiveWindow.Zoom = 100
